$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Lama1"
$ws.Range("C2").Value = "Itgb1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.5587383333333333
$ws.Range("H2").Value = 1.676215
$ws.Range("I2").Value = 0.8486764927018626
$ws.Range("J2").Value = 0.8937587278261895
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 121.928739
$ws.Range("N2").Value = 365.786217
$ws.Range("O2").Value = 0.2282232151508951
$ws.Range("P2").Value = 0.2419720431319445
$ws.Range("Q2").Value = 68.12626041429499
$ws.Range("R2").Value = 613.136343728655
$ws.Range("S2").Value = 0.1936876777874042
$ws.Range("T2").Value = 0.2162646254391106

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Lama1"
$ws.Range("C3").Value = "Itgb1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.5587383333333333
$ws.Range("H3").Value = 1.676215
$ws.Range("I3").Value = 0.8486764927018626
$ws.Range("J3").Value = 0.8937587278261895
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 147.91433
$ws.Range("N3").Value = 443.74299
$ws.Range("O3").Value = 0.2768624053389947
$ws.Range("P3").Value = 0.2935413991166814
$ws.Range("Q3").Value = 82.64540622031666
$ws.Range("R3").Value = 743.8086559828499
$ws.Range("S3").Value = 0.2349666151240995
$ws.Range("T3").Value = 0.262355187438845

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Lama1"
$ws.Range("C4").Value = "Itgb1"
$ws.Range("D4").Value = "Inflammatory-Mac"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.5587383333333333
$ws.Range("H4").Value = 1.676215
$ws.Range("I4").Value = 0.8486764927018626
$ws.Range("J4").Value = 0.8937587278261895
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 83.50496933333334
$ws.Range("N4").Value = 250.514908
$ws.Range("O4").Value = 0.1563025480180701
$ws.Range("P4").Value = 0.1657186665504434
$ws.Range("Q4").Value = 46.65742739035778
$ws.Range("R4").Value = 419.91684651322
$ws.Range("S4").Value = 0.1326502982523402
$ws.Range("T4").Value = 0.1481125045931768

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Lama1"
$ws.Range("C5").Value = "Itgb1"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.5587383333333333
$ws.Range("H5").Value = 1.676215
$ws.Range("I5").Value = 0.8486764927018626
$ws.Range("J5").Value = 0.8937587278261895
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 91.06846250000001
$ws.Range("N5").Value = 182.136925
$ws.Range("O5").Value = 0.1704597085236707
$ws.Range("P5").Value = 0.1204857969594293
$ws.Range("Q5").Value = 50.88344095647917
$ws.Range("R5").Value = 305.3006457388751
$ws.Range("S5").Value = 0.1446651475768507
$ws.Range("T5").Value = 0.1076852326115841

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Lama1"
$ws.Range("C6").Value = "Itgb1"
$ws.Range("D6").Value = "Resolving-Mac"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.5587383333333333
$ws.Range("H6").Value = 1.676215
$ws.Range("I6").Value = 0.8486764927018626
$ws.Range("J6").Value = 0.8937587278261895
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 89.83563
$ws.Range("N6").Value = 269.50689
$ws.Range("O6").Value = 0.1681521229683693
$ws.Range("P6").Value = 0.1782820942415013
$ws.Range("Q6").Value = 50.19461018015
$ws.Range("R6").Value = 451.75149162135
$ws.Range("S6").Value = 0.142706753961168
$ws.Range("T6").Value = 0.159341177743473

# Row 7
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Lama1"
$ws.Range("C7").Value = "Itgb1"
$ws.Range("D7").Value = "ECs"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.09962599999999999
$ws.Range("H7").Value = 0.199252
$ws.Range("I7").Value = 0.1513235072981373
$ws.Range("J7").Value = 0.1062412721738106
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 121.928739
$ws.Range("N7").Value = 365.786217
$ws.Range("O7").Value = 0.2282232151508951
$ws.Range("P7").Value = 0.2419720431319445
$ws.Range("Q7").Value = 12.147272551614
$ws.Range("R7").Value = 72.88363530968398
$ws.Range("S7").Value = 0.03453553736349084
$ws.Range("T7").Value = 0.02570741769283395

# Row 8
$ws.Range("A8").Value = "MuSCs"
$ws.Range("B8").Value = "Lama1"
$ws.Range("C8").Value = "Itgb1"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.09962599999999999
$ws.Range("H8").Value = 0.199252
$ws.Range("I8").Value = 0.1513235072981373
$ws.Range("J8").Value = 0.1062412721738106
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 147.91433
$ws.Range("N8").Value = 443.74299
$ws.Range("O8").Value = 0.2768624053389947
$ws.Range("P8").Value = 0.2935413991166814
$ws.Range("Q8").Value = 14.73611304058
$ws.Range("R8").Value = 88.41667824347999
$ws.Range("S8").Value = 0.04189579021489522
$ws.Range("T8").Value = 0.03118621167783651

# Row 9
$ws.Range("A9").Value = "MuSCs"
$ws.Range("B9").Value = "Lama1"
$ws.Range("C9").Value = "Itgb1"
$ws.Range("D9").Value = "Inflammatory-Mac"
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.09962599999999999
$ws.Range("H9").Value = 0.199252
$ws.Range("I9").Value = 0.1513235072981373
$ws.Range("J9").Value = 0.1062412721738106
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 83.50496933333334
$ws.Range("N9").Value = 250.514908
$ws.Range("O9").Value = 0.1563025480180701
$ws.Range("P9").Value = 0.1657186665504434
$ws.Range("Q9").Value = 8.319266074802666
$ws.Range("R9").Value = 49.915596448816
$ws.Range("S9").Value = 0.02365224976572989
$ws.Range("T9").Value = 0.01760616195726661

# Row 10
$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Lama1"
$ws.Range("C10").Value = "Itgb1"
$ws.Range("D10").Value = "MuSCs"
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.09962599999999999
$ws.Range("H10").Value = 0.199252
$ws.Range("I10").Value = 0.1513235072981373
$ws.Range("J10").Value = 0.1062412721738106
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 91.06846250000001
$ws.Range("N10").Value = 182.136925
$ws.Range("O10").Value = 0.1704597085236707
$ws.Range("P10").Value = 0.1204857969594293
$ws.Range("Q10").Value = 9.072786645025
$ws.Range("R10").Value = 36.2911465801
$ws.Range("S10").Value = 0.02579456094682005
$ws.Range("T10").Value = 0.01280056434784521

# Row 11
$ws.Range("A11").Value = "MuSCs"
$ws.Range("B11").Value = "Lama1"
$ws.Range("C11").Value = "Itgb1"
$ws.Range("D11").Value = "Resolving-Mac"
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.09962599999999999
$ws.Range("H11").Value = 0.199252
$ws.Range("I11").Value = 0.1513235072981373
$ws.Range("J11").Value = 0.1062412721738106
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 89.83563
$ws.Range("N11").Value = 269.50689
$ws.Range("O11").Value = 0.1681521229683693
$ws.Range("P11").Value = 0.1782820942415013
$ws.Range("Q11").Value = 8.949964474379998
$ws.Range("R11").Value = 53.69978684628
$ws.Range("S11").Value = 0.02544536900720132
$ws.Range("T11").Value = 0.01894091649802829
